# This script applies the textual corrections (typo/spelling fixes) that the
# author made to the document. Word's "Editor"/spell-check pass also wraps
# a handful of foreign-language proper nouns (Cybenko, Hornik, relu) in
# proofing marks, but those are internal proofing-engine artifacts that are
# not exposed through the Word object model, so only the content edits below
# are reproducible via COM automation.

$d = $word.ActiveDocument

# 1) "...סכום הנזגרו החלקיות..." -> "...סכום הנגזרות החלקיות..."
$d.Content.Find.Execute(
    "הנזגרו", $true, $false, $false, $false, $false,
    $true, 1, $false, "הנגזרות", 2) | Out-Null

# 2) "...שסכימת ערך פונקצית הסתברות..." -> "...שסכימת ערך פונקציית הסתברות..."
$d.Content.Find.Execute(
    "פונקצית הסתברות", $true, $false, $false, $false, $false,
    $true, 1, $false, "פונקציית הסתברות", 2) | Out-Null

# 3) "בהנתן " -> "בהינתן "
$d.Content.Find.Execute(
    "בהנתן ", $true, $false, $false, $false, $false,
    $true, 1, $false, "בהינתן ", 2) | Out-Null

# 4) "...פונק' התסברות" -> "...פונק' הסתברות"
$d.Content.Find.Execute(
    "התסברות", $true, $false, $false, $false, $false,
    $true, 1, $false, "הסתברות", 2) | Out-Null
